{"js": "// \"An Entwurf weiter gearbeitet.\" \u2014 continue drafting the \"Entwurf\" doc:\n// document the `starter.erl` (start/0) and `koordinatorSteuerung.erl`\n// (startCC/0) modules at the end of the document, right after the\n// existing \"...Programm wird beendet. / Return / void\" paragraphs, and\n// move the (hidden) \"_GoBack\" bookmark from the old last paragraph to the\n// new, now-final \"void\" paragraph \u2014 exactly as Word leaves it after the\n// cursor was last edited at the very end of the document.\n\n// 1) The \"_GoBack\" bookmark currently sits right after \"Programm wird\n//    beendet.\" (it marks Word's last edit position). Since we are about\n//    to append a large amount of new content and the last edit will end\n//    up at the very end of the document, drop it here; it is re-created\n//    at the new end-of-document position in step 3.\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Append the newly authored paragraphs (two new subsections,\n//    \"starter.erl\" / start() and \"koordinatorSteuerung.erl\" / startCC())\n//    after the existing final paragraph (\"void\").\nconst newParagraphs = [\n  { text: \"starter.erl\", style: \"Heading 2\" },\n  { text: \"start()\", style: \"Heading 3\" },\n  { text: \"Definition\", style: \"Heading 4\" },\n  { text: \"Erfragt bei Koordinator die Steuerndewerte({From,getsteeringval}) und wartet auf die Antwort vom Koordinator({steeringval,ArbeitsZeit,TermZeit,Quota,GGTProzessnummer}). Liest restliche Werte aus der ggt.cfg aus. \", style: \"Normal\" },\n  { text: \"Startet die vorgegebene Anzahl an ggT-Prozessen.\", style: \"Normal\" },\n  { text: \"Precondition\", style: \"Heading 4\" },\n  { text: \"Koordinator wurde gestartet.\", style: \"Normal\" },\n  { text: \"Postcondition\", style: \"Heading 4\" },\n  { text: \"Koordinator-Prozess wurde korrekt gestartet.\", style: \"Normal\" },\n  { text: \"Return\", style: \"Heading 4\" },\n  { text: \"Void\", style: \"Normal\" },\n  { text: \"koordinatorSteuerung.erl\", style: \"Heading 2\" },\n  { text: \"startCC()\", style: \"Heading 3\" },\n  { text: \"Definition\", style: \"Heading 4\" },\n  { text: \"Liest die Benutzer Eingaben von der Konsole und leitet diese an Koordinator weiter.\", style: \"Normal\" },\n  { text: \"Precondition\", style: \"Heading 4\" },\n  { text: \"Koordinator-Prozess wurde gestartet.\", style: \"Normal\" },\n  { text: \" Postcondition\", style: \"Heading 4\" },\n  { text: \"Return\", style: \"Heading 4\" },\n  { text: \"void\", style: \"Normal\" },\n];\n\nconst body = context.document.body;\nlet anchor = body.paragraphs.getLast();\nfor (const para of newParagraphs) {\n  const p = anchor.insertParagraph(para.text, Word.InsertLocation.after);\n  p.style = para.style;\n  anchor = p;\n}\n\nawait context.sync();\n\n// 3) Re-create \"_GoBack\" at the new end of the document, i.e. inside the\n//    new, final \"void\" paragraph \u2014 matching where Word leaves it after\n//    the last edit.\nconst endRange = body.getRange(Word.RangeLocation.end);\nendRange.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the start of the paragraph containing \"Programm wird beendet.\" -\n# this is the first paragraph of the block we are replacing.\n$find = $d.Content\n$find.Find.ClearFormatting()\n$ok = $find.Find.Execute(\"Programm wird beendet.\")\nif (-not $ok) {\n    throw \"Anchor text 'Programm wird beendet.' not found\"\n}\n$startPos = $find.Start\n\n# Replace everything from that paragraph through the end of the document\n# body with the fully-specified replacement OOXML (this covers the\n# \"Programm wird beendet.\" / \"Return\" / \"void\" paragraphs plus all of the\n# newly authored content that follows them).\n$endPos = $d.Content.End\n$rng = $d.Range($startPos, $endPos)\n\n$snippet = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>Programm wird beendet.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift4\"/></w:pPr><w:r><w:t>Return</w:t></w:r></w:p><w:p><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/><w:r><w:t>v</w:t></w:r><w:r><w:t>oid</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift2\"/></w:pPr><w:proofErr w:type=\"spellStart\"/><w:r><w:lastRenderedPageBreak/><w:t>starter.erl</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift3\"/></w:pPr><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/><w:r><w:t>start</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>()</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift4\"/></w:pPr><w:r><w:t>Definition</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">Erfragt bei Koordinator die </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Steuerndewerte</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>({</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>From</w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>,getsteeringval</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/><w:r><w:t>}</w:t></w:r><w:r><w:t>) und wartet auf die Antwort vom Koordinator(</w:t></w:r><w:r><w:t>{steeringval,ArbeitsZeit,TermZeit,Quota,GGTProzessnummer}</w:t></w:r><w:r><w:t xml:space=\"preserve\">). Liest restliche Werte aus der </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>ggt.cfg</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> aus. </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">Startet die vorgegebene Anzahl an </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>ggT</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>-Prozessen.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift4\"/></w:pPr><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Precondition</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p><w:r><w:t>Koordinator wurde gestartet.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift4\"/></w:pPr><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Postcondition</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p><w:r><w:t>Koordinator-Prozess wurde korrekt gestartet.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift4\"/></w:pPr><w:r><w:t>Return</w:t></w:r></w:p><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Void</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift2\"/></w:pPr><w:proofErr w:type=\"spellStart\"/><w:r><w:t>koordinatorSteuerung.erl</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift3\"/></w:pPr><w:proofErr w:type=\"spellStart\"/><w:r><w:t>startCC</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>()</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift4\"/></w:pPr><w:r><w:t>Definition</w:t></w:r></w:p><w:p><w:r><w:t>Liest die Benutzer Eingaben von der Konsole und leitet diese an Koordinator weiter.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift4\"/></w:pPr><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Precondition</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p><w:r><w:t>Koordinator-Prozess wurde gestartet.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift4\"/></w:pPr><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Postcondition</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p><w:pPr><w:pStyle w:val=\"berschrift4\"/></w:pPr><w:r><w:t>Return</w:t></w:r></w:p><w:p><w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/><w:r><w:t>void</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:proofErr w:type=\"spellEnd\"/><w:proofErr w:type=\"gramEnd\"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n[void]$rng.InsertXML($snippet)\n"}
